# Fix "missing first row" bug: the original sheet1 had its data rows
# mislabeled (row index 2 was blank/missing because the real data lived
# under row indices 3 and 4). Deleting the phantom blank row 2 shifts the
# real data up so it lines up with its row labels again.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows(2).Delete()

# New columns D/E: first_event / second_event
$ws1.Range("D1").Value = "first_event"
$ws1.Range("E1").Value = "second_event"

$ws1.Range("D2").Value = 40179.958333333
$ws1.Range("D2").NumberFormat = "YYYY-MM-DD\ HH:MM:SS"
$ws1.Range("E2").Value = 41274
$ws1.Range("E2").NumberFormat = "YYYY-MM-DD"

$ws1.Range("D3").Value = 40179
$ws1.Range("D3").NumberFormat = "YYYY-MM-DD\ HH:MM:SS"
$ws1.Range("E3").Value = 41274
$ws1.Range("E3").NumberFormat = "YYYY-MM-DD"

# B3 becomes a formula referencing B2 instead of the literal 5
$ws1.Range("B3").Formula = "=B2"

# New Sheet2, placed after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = 2003
$ws2.Range("B1").NumberFormat = "@"
$ws2.Range("B1").Value = "01"
$ws2.Range("B1").Style = "Normal"
$ws2.Range("C1").Value = 343.12

$ws2.Range("A2").Value = 2003
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "02"
$ws2.Range("B2").Style = "Normal"
$ws2.Range("C2").Value = 345.12

$ws1.Select()
